$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new rows 16-18 with the new default dictionary attributes
$ws.Range("A16").Value = 516
$ws.Range("B16").Value = "CREATION_DATE"
$ws.Range("C16").Value = "mdex:dateTime"
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = "InitDate"

$ws.Range("A17").Value = 516
$ws.Range("B17").Value = "EVENT_DATE"
$ws.Range("C17").Value = "mdex:dateTime"
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = "AwaitP"

$ws.Range("A18").Value = 516
$ws.Range("B18").Value = "NEED_BY_DATE"
$ws.Range("C18").Value = "mdex:dateTime"
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = "Estimated Availability Date"

# Update selection to reflect the new range used
$ws.Range("A2:E18").Select()
